$d = $word.ActiveDocument

# Position an empty, collapsed range at the very end of the document body
# (just before the final section break), then inject the new paragraphs as
# raw WordprocessingML so we get exact control over paragraph/run boundaries.
$endRange = $d.Range($d.Content.End, $d.Content.End)

$w = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$xml = "<w:p $w/>" +
       "<w:p $w>" +
         "<w:r><w:t>Update transakci, paymentů – každ</w:t></w:r>" +
         "<w:r><w:t>é</w:t></w:r>" +
         "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
         "<w:r><w:t>2</w:t></w:r>" +
         "<w:r><w:t xml:space='preserve'> min</w:t></w:r>" +
       "</w:p>" +
       "<w:p $w/>"

$endRange.InsertXML($xml) | Out-Null
